# Add new match-result rows (319-331) to the "Partidos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")
$null = $ws.Activate()

# Each entry: fecha, jugador, equipo, posicion, goles, autogoles, arquero, goles_recibidos, amarillas, rojas, asistencias, penales_atajados
$newRows = @(
    @(45843, "Gember Marin Sarria",      "Azul",     "Arquero",       0, 0, $true,  6, 0, 0, 0, 0),
    @(45843, "Jefferson Delgado",        "Azul",     "Mediocampista", 1, 0, $false, 0, 0, 0, 0, 0),
    @(45843, "Andres Tangarife",         "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @(45843, "Carlos Fernando Valencia", "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 1, 0),
    @(45843, "Cesar Augusto Estrada",    "Azul",     "Delantero",     2, 0, $false, 0, 0, 0, 0, 0),
    @(45843, "Sebastian Giraldo",        "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @(45843, "Armando Vieras",           "Amarillo", "Arquero",       1, 0, $true,  3, 0, 0, 0, 0),
    @(45843, "Invitado",                 "Amarillo", "Arquero",       0, 0, $true,  2, 0, 0, 0, 0),
    @(45843, "Andres Jurado",            "Amarillo", "Delantero",     1, 0, $false, 0, 0, 0, 1, 0),
    @(45843, "Hermes Marquez",           "Amarillo", "Defensa",       1, 0, $false, 0, 0, 0, 0, 0),
    @(45843, "Edwing Yesid Castillo",    "Amarillo", "Mediocampista", 1, 0, $false, 0, 0, 0, 1, 0),
    @(45843, "Julian Mbappe",            "Amarillo", "Mediocampista", 2, 0, $false, 0, 0, 0, 0, 0),
    @(45843, "Carlos Julio Delgado",     "Amarillo", "Defensa",       0, 0, $false, 0, 0, 0, 1, 0)
)

$startRow = 319
$endRow = $startRow + $newRows.Count - 1

# Copy the date-cell formatting (style index already used by column A) onto
# the new A319:A331 cells before assigning values, so the existing style is
# reused instead of a new (duplicate) one being created.
$null = $ws.Cells.Item($startRow - 1, 1).Copy()
$null = $ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
    $ws.Cells.Item($r, 12).Value = $data[11]
}

# Update selection/active cell to mirror the final edited state
$null = $ws.Range("B333").Select()
